# Sprint 39 - Day 4 test case summary: record counts for the purchase page
# test cases that were created and executed (Total Written / Total
# Execution / Total Review for the "Day 4" block, rows 25-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

$ws.Range("C25").Value = 884
$ws.Range("C26").Value = 1115
$ws.Range("C27").Value = 636

# Re-create the merged header cells so their stored order matches a fresh
# Excel save (Excel re-emits <mergeCells> sorted by sheet order after an
# edit like this).
$mergedRanges = @(
    "B2:C2", "B8:C8", "B16:C16", "B24:C24", "B30:C30",
    "B36:C36", "B42:C42", "B48:C48", "B54:C54", "B61:C61"
)
foreach ($r in $mergedRanges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $mergedRanges) {
    $ws.Range($r).Merge()
}

# Leave the cursor on the last cell that was touched, like a user would
# after typing the Day 4 review count.
$ws.Range("C27").Select()
